# DE223800: upload templates missing mandatory asterisk
# Prefix the header labels on the "Disty to Direct" sheet with "*" to
# indicate that the fields are mandatory.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Disty to Direct")

$ws.Range("A5").Value = "*Group ID"
$ws.Range("B5").Value = "*Node Type"
$ws.Range("C5").Value = "*Sales Finance Hierarchy"
$ws.Range("D5").Value = "*Node Code"
